# home page method updated
# Updates execution timestamps for all test rows and flips the
# "resend otp" block test result from PASSED to FAILED, while the
# "Home Page Loads Successfully" row picks up the final timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout: A=Test Case ID, B=Test Case Name, C=Status, D=Execution Time, E=Comment

# Row 2 - User Login with Valid Credentials
$ws.Range("D2").Value = "21/04/2025 11:03:00 AM"

# Row 3 - Verify empty email state
$ws.Range("D3").Value = "21/04/2025 11:03:05 AM"

# Row 4 - Verify Login with Invalid Email
$ws.Range("D4").Value = "21/04/2025 11:03:06 AM"

# Row 5 - Verify Login with Not Registred Email
$ws.Range("D5").Value = "21/04/2025 11:03:07 AM"

# Row 6 - Verify admin viewer Shouldn't be able to login
$ws.Range("D6").Value = "21/04/2025 11:03:07 AM"

# Row 7 - Verify Login with Wrong OTP
$ws.Range("D7").Value = "21/04/2025 11:03:10 AM"

# Row 8 - Verify same email on OTP page
$ws.Range("D8").Value = "21/04/2025 11:03:10 AM"

# Row 9 - Verify account block after attempting wrong OTP for 5 times
$ws.Range("D9").Value = "21/04/2025 11:03:23 AM"

# Row 10 - Verify Go To Sign In page Navigation
$ws.Range("D10").Value = "21/04/2025 11:03:25 AM"

# Row 11 - Verify that navigation and getOTP blocked for blocked account
$ws.Range("D11").Value = "21/04/2025 11:03:26 AM"

# Row 12 - Verify Resend OTP button
$ws.Range("D12").Value = "21/04/2025 11:04:30 AM"

# Row 13 - Verify account block after 5 times of resend otp
# Test now fails with a mismatched error message.
$ws.Range("C13").Value = "FAILED"
$ws.Range("D13").Value = "21/04/2025 11:07:42 AM"
$ws.Range("E13").Value = "Error message isn't as expected expected [You have reached the maximum login attempts for the day. Please try again after 24 hours.] but found [Max OTP retry limit reached. Please try again later]"

# Row 14 - Verify Home Page Loads Successfully
$ws.Range("D14").Value = "21/04/2025 11:07:42 AM"
